$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 4
$ws.Range("A4").Value = 110282846
$ws.Range("B4").Value = 103288
$ws.Range("E4").Value = 221144
$ws.Range("I4").Value = "'10"
$ws.Range("I4").Style = "Normal"
$ws.Range("D4").Value = "LC"
$ws.Range("F4").Value = "Grönpyrola"
$ws.Range("G4").Value = "Pyrola chlorantha"
$ws.Range("H4").Value = "Sw."
$ws.Range("K4").Value = "blomning"

# Row 5
$ws.Range("A5").Value = 110282835
$ws.Range("B5").Value = 89405
$ws.Range("E5").Value = 1202
$ws.Range("Q5").Value = 600805.3583702671
$ws.Range("R5").Value = 6613969.910894822
$ws.Range("I5").Value = "'1"
$ws.Range("I5").Style = "Normal"
$ws.Range("D5").Value = "NT"
$ws.Range("F5").Value = "Ullticka"
$ws.Range("G5").Value = "Phellinidium ferrugineofuscum"
$ws.Range("H5").Value = "(P.Karst.) Fiasson & Niemelä"
$ws.Range("J5").Value = "mycel"
$ws.Range("K5").ClearContents()
$ws.Range("N4").Copy($ws.Range("K5"))
$ws.Range("L5").ClearContents()
$ws.Range("AF5").ClearContents()

# Row 6
$ws.Range("A6").Value = 110282764
$ws.Range("B6").Value = 96348
$ws.Range("E6").Value = 220787
$ws.Range("Q6").Value = 600749.0751519018
$ws.Range("R6").Value = 6613971.934424319
$ws.Range("I6").Value = "'10"
$ws.Range("I6").Style = "Normal"
$ws.Range("D6").Value = "VU"
$ws.Range("F6").Value = "Knärot"
$ws.Range("G6").Value = "Goodyera repens"
$ws.Range("H6").Value = "(L.) R. Br."
$ws.Range("J6").Value = "plantor/tuvor"
$ws.Range("K6").Value = "fullt utvecklade blad"
$ws.Range("L6").ClearContents()
$ws.Range("N4").Copy($ws.Range("L6"))
$ws.Range("AF6").ClearContents()
$ws.Range("N4").Copy($ws.Range("AF6"))

# Row 7
$ws.Range("A7").Value = 110282856
$ws.Range("B7").Value = 89802
$ws.Range("E7").Value = 5420
$ws.Range("Q7").Value = 600677.6983460309
$ws.Range("R7").Value = 6613951.301940188
$ws.Range("F7").Value = "Grovticka"
$ws.Range("G7").Value = "Phaeolus schweinitzii"
$ws.Range("H7").Value = "(Fr.) Pat."
$ws.Range("J7").Value = "fruktkroppar"

# Row 8
$ws.Range("A8").Value = 110282836
$ws.Range("B8").Value = 89793
$ws.Range("E8").Value = 4217
$ws.Range("D8").Value = "LC"
$ws.Range("F8").Value = "Blodticka"
$ws.Range("G8").Value = "Meruliopsis taxicola"
$ws.Range("H8").Value = "(Pers.:Fr.) Bondartsev"

# Row 9
$ws.Range("A9").Value = 110282848
$ws.Range("B9").Value = 96348
$ws.Range("E9").Value = 220787
$ws.Range("I9").Value = "'5"
$ws.Range("I9").Style = "Normal"
$ws.Range("D9").Value = "VU"
$ws.Range("F9").Value = "Knärot"
$ws.Range("G9").Value = "Goodyera repens"
$ws.Range("H9").Value = "(L.) R. Br."
$ws.Range("K9").Value = "fullt utvecklade blad"
